# Fruta / hortaliza, semanal
# Inserts two new daily price rows (44-45) for "Cereza" variedad "Santina"
# recorded on serial date 44587 (2022-01-26), pushing the previously
# existing rows 44-98 down to 46-100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 44-45; this shifts old rows 44-98 down to 46-100
# (all their cell values, including the D-column date style, move with them).
$ws.Rows("44:45").Insert()

# --- Row 44: Cereza / Santina / Primera ---
$ws.Cells.Item(44, 1).Value = 7
$ws.Cells.Item(44, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value = "Ñuble"
$ws.Cells.Item(44, 4).Value = 44587
$ws.Cells.Item(44, 5).Value = 16
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100103
$ws.Cells.Item(44, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(44, 9).Value = 100103001
$ws.Cells.Item(44, 10).Value = "Cereza"
$ws.Cells.Item(44, 11).Value = "Santina"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 400
$ws.Cells.Item(44, 14).Value = 5500
$ws.Cells.Item(44, 15).Value = 6000
$ws.Cells.Item(44, 16).Value = 5750
$ws.Cells.Item(44, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(44, 19).Value = 575
$ws.Cells.Item(44, 20).Value = 10

# --- Row 45: Cereza / Santina / Segunda ---
$ws.Cells.Item(45, 1).Value = 7
$ws.Cells.Item(45, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(45, 3).Value = "Ñuble"
$ws.Cells.Item(45, 4).Value = 44587
$ws.Cells.Item(45, 5).Value = 16
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100103
$ws.Cells.Item(45, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(45, 9).Value = 100103001
$ws.Cells.Item(45, 10).Value = "Cereza"
$ws.Cells.Item(45, 11).Value = "Santina"
$ws.Cells.Item(45, 12).Value = "Segunda"
$ws.Cells.Item(45, 13).Value = 200
$ws.Cells.Item(45, 14).Value = 4500
$ws.Cells.Item(45, 15).Value = 5000
$ws.Cells.Item(45, 16).Value = 4750
$ws.Cells.Item(45, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(45, 19).Value = 475
$ws.Cells.Item(45, 20).Value = 10

# Ensure the date cells keep the workbook's date/time number format
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
